$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M18").Value = -939.9000000000001
$ws.Range("I18").Value = 1223.9
$ws.Range("K18").Value = 1223.9
$ws.Range("H18").Value = 1353.25
$ws.Range("J112").Value = 2932.9333
$ws.Range("L112").Value = 8798.7999
$ws.Range("K112").Value = 2100
$ws.Range("M112").Value = -992
$ws.Range("H112").Value = 2462.842
$ws.Range("N112").Value = -11014.7999
$ws.Range("I112").Value = 700
$ws.Range("M116").Value = 1905.9333
$ws.Range("I116").Value = 1536.0667
$ws.Range("H116").Value = 2315.2273
$ws.Range("K116").Value = 1536.0667
$ws.Range("N121").Value = -6389
$ws.Range("L121").Value = 2895
$ws.Range("H121").Value = 965
$ws.Range("J121").Value = 965
$ws.Range("N133").Value = -43108.332
$ws.Range("L133").Value = 32988.332
$ws.Range("H133").Value = 32988.332
$ws.Range("J133").Value = 32988.332
$ws.Range("K137").Value = 3329.1666
$ws.Range("J137").Value = 1800.3334
$ws.Range("M137").Value = -779.1665999999996
$ws.Range("L137").Value = 5401.0002
$ws.Range("H137").Value = 1455.0278
$ws.Range("I137").Value = 1109.7222
$ws.Range("N137").Value = -10501.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I74").Value = 985.7059
$ws.Range("J74").Value = 3960
$ws.Range("L74").Value = 3960
$ws.Range("K74").Value = 985.7059
$ws.Range("H74").Value = 1661.6818
$ws.Range("N74").Value = -5708
$ws.Range("M74").Value = -111.7059
$ws.Range("H77").Value = 1661.6818
$ws.Range("K77").Value = 4928.529500000001
$ws.Range("I77").Value = 985.7059
$ws.Range("J77").Value = 3960
$ws.Range("L77").Value = 19800
$ws.Range("M77").Value = -560.5295000000006
$ws.Range("N77").Value = -28536
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("J135").Value = 0
$ws.Range("H135").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M105").Value = -62499743
$ws.Range("K105").Value = 62501490
$ws.Range("H105").Value = 50001560
$ws.Range("I105").Value = 62501490
$ws.Range("N116").Value = -79177.664
$ws.Range("J116").Value = 69999.664
$ws.Range("H116").Value = 69999.664
$ws.Range("L116").Value = 69999.664
$ws.Range("N132").Value = -44200
$ws.Range("L132").Value = 34080
$ws.Range("J132").Value = 34080
$ws.Range("H132").Value = 34080

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I17").Value = 508.5
$ws.Range("K17").Value = 508.5
$ws.Range("M17").Value = -334.5
$ws.Range("H17").Value = 508.5
$ws.Range("I31").Value = 1418.0541
$ws.Range("H31").Value = 1511.7
$ws.Range("M31").Value = -1123.0541
$ws.Range("K31").Value = 1418.0541
$ws.Range("K34").Value = 1418.0541
$ws.Range("I34").Value = 1418.0541
$ws.Range("H34").Value = 1511.7
$ws.Range("M34").Value = -1216.0541
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H41").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("J50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H51").Value = 23400
$ws.Range("N51").Value = -25722
$ws.Range("J51").Value = 24250
$ws.Range("L51").Value = 24250
$ws.Range("M59").Value = 695
$ws.Range("K59").Value = 450
$ws.Range("I59").Value = 450
$ws.Range("H59").Value = 24075
$ws.Range("N59").Value = -31090
$ws.Range("J59").Value = 28800
$ws.Range("L59").Value = 28800
$ws.Range("K60").Value = 6990
$ws.Range("I60").Value = 6990
$ws.Range("H60").Value = 6990
$ws.Range("M60").Value = -6479
$ws.Range("N61").Value = -24946
$ws.Range("J61").Value = 24250
$ws.Range("H61").Value = 23400
$ws.Range("L61").Value = 24250
$ws.Range("I74").Value = 14000
$ws.Range("K74").Value = 14000
$ws.Range("H74").Value = 23500
$ws.Range("M74").Value = -13126
$ws.Range("H75").Value = 11573.333
$ws.Range("N75").Value = -13569.333
$ws.Range("J75").Value = 11573.333
$ws.Range("L75").Value = 11573.333
$ws.Range("H77").Value = 23500
$ws.Range("K77").Value = 42000
$ws.Range("I77").Value = 14000
$ws.Range("M77").Value = -37632
$ws.Range("H78").Value = 11573.333
$ws.Range("J78").Value = 11573.333
$ws.Range("N78").Value = -44703.999
$ws.Range("L78").Value = 34719.999
$ws.Range("K99").Value = 0
$ws.Range("H99").Value = 1800
$ws.Range("I99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("I122").Value = 2400
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -4750
$ws.Range("K122").Value = 7200
$ws.Range("H122").Value = 1600
$ws.Range("N122").Value = -8500
$ws.Range("J122").Value = 1200
$ws.Range("M126").ClearContents()
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H126").Value = 1800
$ws.Range("I132").Value = 956.2632
$ws.Range("N132").Value = -11260.1819
$ws.Range("L132").Value = 6200.1819
$ws.Range("K132").Value = 2868.7896
$ws.Range("J132").Value = 2066.7273
$ws.Range("H132").Value = 1363.4333
$ws.Range("M132").Value = -338.7896000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 1700
$ws.Range("N87").Value = -11496
$ws.Range("L87").Value = 9000
$ws.Range("J87").Value = 3000
$ws.Range("H90").Value = 1700
$ws.Range("J90").Value = 3000
$ws.Range("N90").Value = -39480
$ws.Range("L90").Value = 27000
$ws.Range("K107").Value = 1362
$ws.Range("H107").Value = 5876.5557
$ws.Range("L107").Value = 23886.462
$ws.Range("J107").Value = 7962.154
$ws.Range("I107").Value = 454
$ws.Range("N107").Value = -27726.462
$ws.Range("M107").Value = 558
$ws.Range("L129").Value = 13161847.5
$ws.Range("M129").Value = -166664008
$ws.Range("K129").Value = 166669008
$ws.Range("J129").Value = 4387282.5
$ws.Range("N129").Value = -13171847.5
$ws.Range("I129").Value = 55556336
$ws.Range("H129").Value = 16667855
$ws.Range("H131").Value = 12821382
$ws.Range("J131").Value = 976.65576
$ws.Range("L131").Value = 2929.96728
$ws.Range("N131").Value = -13009.96728

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J123").Value = 21114.143
$ws.Range("H123").Value = 21114.143
$ws.Range("L123").Value = 21114.143
$ws.Range("N123").Value = -26014.143

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M40").Value = -2779.6
$ws.Range("K40").Value = 2915.6
$ws.Range("H40").Value = 4582.8
$ws.Range("N40").Value = -6522
$ws.Range("J40").Value = 6250
$ws.Range("L40").Value = 6250
$ws.Range("I40").Value = 2915.6
$ws.Range("I122").Value = 35716376
$ws.Range("L122").Value = 7507.5
$ws.Range("M122").Value = -107146678
$ws.Range("K122").Value = 107149128
$ws.Range("H122").Value = 27779960
$ws.Range("N122").Value = -12407.5
$ws.Range("J122").Value = 2502.5
$ws.Range("J123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M136").Value = -658.7586000000001
$ws.Range("K136").Value = 3208.7586
$ws.Range("J136").Value = 1557
$ws.Range("N136").Value = -9771
$ws.Range("H136").Value = 1212.2439
$ws.Range("L136").Value = 4671
$ws.Range("I136").Value = 1069.5862
